$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 date/time
$ws.Range("A2").Value = "2025-04-08T11:13"

# Update C2:F2 from numeric values to text values (force text with a quote prefix,
# like a user typing '1 / '1,1 / '1,88 into the cell, so they stay text not numbers)
$ws.Range("C2").Value = "'1"
$ws.Range("D2").Value = "'1"
$ws.Range("E2").Value = "'1,1"
$ws.Range("F2").Value = "'1,88"

# Update Q2:S2 from "t" to "0"
$ws.Range("Q2").Value = "'0"
$ws.Range("R2").Value = "'0"
$ws.Range("S2").Value = "'0"

# Delete entire row 3 (shifts nothing up since it's the last row, removes it)
$ws.Rows("3").Delete()
